$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '305.03'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '1.34%'
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '35.89'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '1.18%'
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.032'
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '-0.67%'
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.08096'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '1.11%'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.918'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '1.30%'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '4.141'
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '2.37%'
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '7.846'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.9320'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '0.57%'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.1283'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '-12.82%'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.1917'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '0.87%'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.09194'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '2.21%'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.03490'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '1.75%'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.09929'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '0.85%'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.001415'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '1.09%'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.006634'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '13.95%'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '2.20%'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.056'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '4.15%'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.3447'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '0.10%'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.183'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '2.46%'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '0.11%'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.2532'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '5.66%'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.04406'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '-1.82%'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.001236'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '2.03%'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.004726'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '-0.79%'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0001301'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '5.87%'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0003131'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '3.65%'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01980'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '7.90%'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.05196'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '9.02%'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '2.48%'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.01010'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '-4.89%'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1372'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '3.37%'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.002102'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '-0.30%'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.01069'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '-1.76%'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00006355'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '2.06%'
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000751'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '0.19%'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '64.96'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '0.45%'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.001601'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '-3.49%'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.00002102'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '0.19%'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0002002'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '0.19%'
$c.Style = "Normal"
